$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 0.0012
$ws.Range("D4").Value = 508.5162
$ws.Range("E4").Value = 16.7674
$ws.Range("F4").Value = 0.55
$ws.Range("G4").Value = 0.07530000000000001
$ws.Range("H4").Value = 0.5002
$ws.Range("I4").Value = 0.0477
$ws.Range("C5").Value = 0.0008
$ws.Range("D5").Value = 478.4193
$ws.Range("E5").Value = 3.9519
$ws.Range("F5").Value = 0.5511
$ws.Range("G5").Value = 0.007900000000000001
$ws.Range("H5").Value = 0.4654
$ws.Range("I5").Value = 0
$ws.Range("C6").Value = 0.0146
$ws.Range("D6").Value = 0.6108
$ws.Range("E6").Value = -25.9108
$ws.Range("F6").Value = 0.3197
$ws.Range("G6").Value = 0.0022
$ws.Range("H6").Value = 0.4817
$ws.Range("I6").Value = 0.0009
$ws.Range("C7").Value = 0.008200000000000001
$ws.Range("D7").Value = 3.1412
$ws.Range("E7").Value = 7.408
$ws.Range("F7").Value = 0.5508999999999999
$ws.Range("G7").Value = 0.0004
$ws.Range("H7").Value = 0.3747
$ws.Range("I7").Value = 0.001
$ws.Range("C8").Value = 0.0168
$ws.Range("D8").Value = 148.9547
$ws.Range("E8").Value = 2.2989
$ws.Range("F8").Value = 86021.67389999999
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 68590.1407
$ws.Range("I8").Value = 0.0005999999999999999
$ws.Range("C9").Value = 0.0115
$ws.Range("D9").Value = 166.5822
$ws.Range("E9").Value = 0.1339
$ws.Range("F9").Value = 114613.889
$ws.Range("G9").Value = 0.004
$ws.Range("H9").Value = 113051.0417
$ws.Range("I9").Value = 0.0013
$ws.Range("C10").Value = 0.0045
$ws.Range("D10").Value = 53.0471
$ws.Range("E10").Value = -0.1748
$ws.Range("F10").Value = 2.089
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 2.1426
$ws.Range("I10").Value = 0.0239
$ws.Range("C11").Value = 0.0437
$ws.Range("D11").Value = 601.4115
$ws.Range("E11").Value = 1.3951
$ws.Range("F11").Value = 6283.6503
$ws.Range("G11").Value = 0.0031
$ws.Range("H11").Value = 4444.9265
$ws.Range("I11").Value = 0.0031
$ws.Range("C12").Value = 0.0245
$ws.Range("D12").Value = 29.4255
$ws.Range("E12").Value = 2.6811
$ws.Range("F12").Value = 1.594
$ws.Range("G12").Value = 0.2259
$ws.Range("H12").Value = 1.515
$ws.Range("I12").Value = 0.06270000000000001
$ws.Range("C13").Value = 0.0091
$ws.Range("D13").Value = 8.411199999999999
$ws.Range("E13").Value = 0.5188
$ws.Range("F13").Value = 0.7977
$ws.Range("G13").Value = 0.2831
$ws.Range("H13").Value = 0.7688
$ws.Range("I13").Value = 0.1621
$ws.Range("D14").Value = 4154.2697
$ws.Range("E14").Value = 13.0684
$ws.Range("F14").Value = 0.2559
$ws.Range("G14").Value = 0.2025
$ws.Range("H14").Value = 0.2794
$ws.Range("I14").Value = 0.2574
$ws.Range("C15").Value = 0.0718
$ws.Range("D15").Value = 1486.0761
$ws.Range("E15").Value = 16.1801
$ws.Range("F15").Value = 0.9927
$ws.Range("G15").Value = 0.4948
$ws.Range("H15").Value = 0.9388
$ws.Range("I15").Value = 0.6129
$ws.Range("C16").Value = 0.0391
$ws.Range("D16").Value = 68.9068
$ws.Range("E16").Value = 0.2887
$ws.Range("F16").Value = 2.5567
$ws.Range("G16").Value = 0.1512
$ws.Range("H16").Value = 2.7311
$ws.Range("I16").Value = 0.0882
$ws.Range("C17").Value = 0.0049
$ws.Range("D17").Value = -0.8939
$ws.Range("E17").Value = 3.778
$ws.Range("F17").Value = 0.5053
$ws.Range("G17").Value = 0.118
$ws.Range("H17").Value = 0.4864
$ws.Range("I17").Value = 0.2481
$ws.Range("C18").Value = 0.0009
$ws.Range("D18").Value = 7.6487
$ws.Range("E18").Value = 1.7074
$ws.Range("F18").Value = 0.482
$ws.Range("G18").Value = 0.0366
$ws.Range("H18").Value = 0.4962
$ws.Range("I18").Value = 0
$ws.Range("D19").Value = 3.0299
$ws.Range("E19").Value = 1.2246
$ws.Range("F19").Value = 2.6395
$ws.Range("G19").Value = 0.0182
$ws.Range("H19").Value = 2.513
$ws.Range("I19").Value = 0.009299999999999999
$ws.Range("C20").Value = 0.0546
$ws.Range("D20").Value = 0.4447
$ws.Range("E20").Value = 5.7306
$ws.Range("F20").Value = 3.0314
$ws.Range("G20").Value = 0.1319
$ws.Range("H20").Value = 2.7474
$ws.Range("I20").Value = 0.1357
$ws.Range("C21").Value = 0.8934
$ws.Range("D21").Value = 507.1132
$ws.Range("E21").Value = -5.181
$ws.Range("F21").Value = 87205.6191
$ws.Range("G21").Value = 0.1832
$ws.Range("H21").Value = 71983.50599999999
$ws.Range("I21").Value = 0.2457
$ws.Range("C22").Value = 0.0303
$ws.Range("D22").Value = -2.3016
$ws.Range("E22").Value = -7.5272
$ws.Range("F22").Value = 3.2474
$ws.Range("G22").Value = 0.0516
$ws.Range("H22").Value = 2.7512
$ws.Range("I22").Value = 0.0679
$ws.Range("C23").Value = 0.6404
$ws.Range("D23").Value = 458.5277
$ws.Range("E23").Value = 4.3641
$ws.Range("F23").Value = 92676.80409999999
$ws.Range("G23").Value = 0.0868
$ws.Range("H23").Value = 76425.2227
$ws.Range("I23").Value = 0.1492
